$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Find the (first remaining) occurrence of $oldText anywhere in the document,
# replace its text with $newText, and "pin" the run by perturbing Font.Color
# so that the engine does not silently re-merge it with a neighbouring run
# that happens to share identical formatting.  Returns the (collapsed-to-
# new-text) Range so callers can anchor further inserts after it.
function Set-RunText($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $true, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output ("NOT FOUND: " + $oldText)
        return $null
    }
    $rng.Text = $newText
    $rng.Font.Color = 255
    $rng.Font.Color = 0
    return $rng
}

# Insert a brand-new run of text immediately after $rng, copying the
# character formatting of $rng, and pin it as its own run.
function Insert-RunAfter($rng, $newText) {
    $fontName = $rng.Font.NameAscii
    $fontSize = $rng.Font.Size
    $newRng = $rng.Duplicate
    $newRng.Collapse(0)   # wdCollapseEnd
    $newRng.InsertAfter($newText)
    $newRng.Font.NameAscii = $fontName
    $newRng.Font.Name = $fontName
    $newRng.Font.Size = $fontSize
    $newRng.Font.Color = 255
    $newRng.Font.Color = 0
    return $newRng
}

# ---------------------------------------------------------------------------
# Title / byline / contact block
# ---------------------------------------------------------------------------

Set-RunText "Universal Constants: Guiding Forces of the Cosmos" "Government - A Guiding Force in Society's Progress" | Out-Null
Set-RunText "Cosmos Explorer" "Elpidio Hernandez" | Out-Null
Set-RunText "cosmic_explorer@astronomicalsociety" "elpidio" | Out-Null

$r = Set-RunText "org" "hernandez765@gmail"
$r = Insert-RunAfter $r "."
$r = Insert-RunAfter $r "com"

# ---------------------------------------------------------------------------
# Body paragraph (first block)
# ---------------------------------------------------------------------------

Set-RunText "Across the vast expanse of the universe, a set of immutable principles governs the fundamental workings of matter, energy, and spacetime" "Governments, in their intricate web of processes and institutions, hold the power to shape the destinies of nations and their citizens" | Out-Null
Set-RunText " Known as universal constants, these values underpin our understanding of the universe and guide our exploration into its mysteries" " From the intricate machinery of global superpowers to the intimate governance of local communities, governments are a testament to human society's endless quest for order, justice, and progress" | Out-Null
Set-RunText " From the speed of light to the mass of the electron, each constant serves as a cornerstone of physics, chemistry, and cosmology" " Their study is an exploration into the art of governing and the science of human nature, an understanding of how leaders guide, policies shape, and institutions underpin the harmonious functioning of societies" | Out-Null

Set-RunText "Universal constants hold the key to understanding the interactions between particles, the properties of materials, and the evolution of the universe" "The existence of governments is intertwined with the complexity of human interactions and the challenges of communal living" | Out-Null
Set-RunText " From the minutest subatomic realm to the grandest cosmic phenomena, the laws of nature are dictated by these fundamental values" " To comprehend their purpose, we embark on a journey through history, examining the birth of civilizations, the evolution of governance structures, and the struggles for power that have shaped our world" | Out-Null

$r = Set-RunText " Like notes within a symphony, universal constants orchestrate the cosmic harmony, providing a framework that allows us to comprehend the intricate mechanisms of the universe" " We dissect the nature of power, its forms, its sources, and its uses, recognizing the inherent responsibility that comes with its possession"
$r = Insert-RunAfter $r "."
$r = Insert-RunAfter $r " We delve into the ideas of democracy, authoritarianism, and the intricate balance between liberty and authority"

Set-RunText "In our quest to unravel the secrets of the cosmos, universal constants serve as invaluable tools" "Governments are more than just abstract concepts; they impact our everyday lives in tangible ways" | Out-Null
Set-RunText " By precisely measuring and analyzing these values, physicists and cosmologists probe the deepest levels of reality, uncovering the underlying principles that govern the universe" " Through the provision of essential services, regulating economic activities, and maintaining law and order, governments create a foundation for stability and prosperity" | Out-Null

$r = Set-RunText " They act as checkpoints against theoretical models, helping us refine our understanding of fundamental laws and pushing the boundaries of scientific knowledge" " We examine the role of taxation, the intricacies of public finance, and the distribution of resources, acknowledging the moral quandaries that arise in allocating societal resources"
$r = Insert-RunAfter $r "."
$r = Insert-RunAfter $r " We analyze the interplay between government and the economy, exploring how policies shape markets, promote growth, and address economic disparities"

# ---------------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------------

Set-RunText "Universal constants are the immutable values that govern the fundamental workings of the universe" "In summary, the study of government is an exploration of power, its forms, and its uses" | Out-Null

$r = Set-RunText " They are essential to understanding the interactions between particles, the properties of materials, and the evolution of the universe" " We examine the historical and theoretical foundations of government, analyzing the nature of leadership, the types of political systems, and the complex relationship between government "
$r = Insert-RunAfter $r "and the people"

Set-RunText " By measuring and analyzing these constants, scientists gain insights into the deepest levels of reality, pushing the boundaries of scientific knowledge" " We investigate the functions of government, including the provision of services, regulation of the economy, and maintenance of law and order" | Out-Null
Set-RunText " They serve as guiding forces, helping us unravel the intricate mechanisms of the cosmos and forge a deeper connection with the mysteries that lie beyond" " Finally, we consider the challenges and opportunities facing governments in the 21st century, emphasizing the need for effective governance in addressing pressing global issues" | Out-Null

# ---------------------------------------------------------------------------
# Trailing empty paragraph added at the end of the document body
# ---------------------------------------------------------------------------

$endRng = $d.Paragraphs.Last.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter() | Out-Null

Write-Output "edit complete"
